# Update countries & provincias Spain
# Applies updated COVID country statistics and re-sorted rank swaps
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Estados Unidos
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 1350194
$ws.Cells.Item(4, 3).Value = 2885
$ws.Cells.Item(4, 4).Value = 238144
$ws.Cells.Item(4, 5).Value = 1031929
$ws.Cells.Item(4, 6).Value = 16816
$ws.Cells.Item(4, 7).Value = 84
$ws.Cells.Item(4, 8).Value = 80121

# Row 11: Brasil
$ws.Cells.Item(11, 1).Value = "Brasil"
$ws.Cells.Item(11, 2).Value = 156604
$ws.Cells.Item(11, 3).Value = 543
$ws.Cells.Item(11, 4).Value = 61685
$ws.Cells.Item(11, 5).Value = 84222
$ws.Cells.Item(11, 6).Value = 8318
$ws.Cells.Item(11, 7).Value = 41
$ws.Cells.Item(11, 8).Value = 10697

# Row 16: India
$ws.Cells.Item(16, 1).Value = "India"
$ws.Cells.Item(16, 2).Value = 65021
$ws.Cells.Item(16, 3).Value = 2213
$ws.Cells.Item(16, 4).Value = 19665
$ws.Cells.Item(16, 5).Value = 43203
$ws.Cells.Item(16, 6).Value = 0
$ws.Cells.Item(16, 7).Value = 52
$ws.Cells.Item(16, 8).Value = 2153

# Row 17: Peru
$ws.Cells.Item(17, 1).Value = "Peru"
$ws.Cells.Item(17, 2).Value = 65015
$ws.Cells.Item(17, 3).Value = 0
$ws.Cells.Item(17, 4).Value = 20246
$ws.Cells.Item(17, 5).Value = 42955
$ws.Cells.Item(17, 6).Value = 748
$ws.Cells.Item(17, 7).Value = 0
$ws.Cells.Item(17, 8).Value = 1814

# Row 58: Argelia
$ws.Cells.Item(58, 1).Value = "Argelia"
$ws.Cells.Item(58, 2).Value = 5723
$ws.Cells.Item(58, 3).Value = 165
$ws.Cells.Item(58, 4).Value = 2678
$ws.Cells.Item(58, 5).Value = 2543
$ws.Cells.Item(58, 6).Value = 22
$ws.Cells.Item(58, 7).Value = 8
$ws.Cells.Item(58, 8).Value = 502

# Row 70: Irak
$ws.Cells.Item(70, 1).Value = "Irak"
$ws.Cells.Item(70, 2).Value = 2767
$ws.Cells.Item(70, 3).Value = 88
$ws.Cells.Item(70, 4).Value = 1734
$ws.Cells.Item(70, 5).Value = 924
$ws.Cells.Item(70, 6).Value = 0
$ws.Cells.Item(70, 7).Value = 2
$ws.Cells.Item(70, 8).Value = 109

# Row 71: Grecia
$ws.Cells.Item(71, 1).Value = "Grecia"
$ws.Cells.Item(71, 2).Value = 2716
$ws.Cells.Item(71, 3).Value = 6
$ws.Cells.Item(71, 4).Value = 1374
$ws.Cells.Item(71, 5).Value = 1191
$ws.Cells.Item(71, 6).Value = 30
$ws.Cells.Item(71, 7).Value = 0
$ws.Cells.Item(71, 8).Value = 151

# Row 96: Mayotte
$ws.Cells.Item(96, 1).Value = "Mayotte"
$ws.Cells.Item(96, 2).Value = 1023
$ws.Cells.Item(96, 3).Value = 35
$ws.Cells.Item(96, 4).Value = 492
$ws.Cells.Item(96, 5).Value = 520
$ws.Cells.Item(96, 6).Value = 9
$ws.Cells.Item(96, 7).Value = 0
$ws.Cells.Item(96, 8).Value = 11

# Row 97: Kirguistan
$ws.Cells.Item(97, 1).Value = "Kirguistan"
$ws.Cells.Item(97, 2).Value = 1002
$ws.Cells.Item(97, 3).Value = 71
$ws.Cells.Item(97, 4).Value = 675
$ws.Cells.Item(97, 5).Value = 315
$ws.Cells.Item(97, 6).Value = 13
$ws.Cells.Item(97, 7).Value = 0
$ws.Cells.Item(97, 8).Value = 12

# Row 101: Republica de Chipre
$ws.Cells.Item(101, 1).Value = "Republica de Chipre"
$ws.Cells.Item(101, 2).Value = 898
$ws.Cells.Item(101, 3).Value = 6
$ws.Cells.Item(101, 4).Value = 401
$ws.Cells.Item(101, 5).Value = 482
$ws.Cells.Item(101, 6).Value = 10
$ws.Cells.Item(101, 7).Value = 0
$ws.Cells.Item(101, 8).Value = 15

# Row 118: Georgia
$ws.Cells.Item(118, 1).Value = "Georgia"
$ws.Cells.Item(118, 2).Value = 635
$ws.Cells.Item(118, 3).Value = 9
$ws.Cells.Item(118, 4).Value = 309
$ws.Cells.Item(118, 5).Value = 316
$ws.Cells.Item(118, 6).Value = 6
$ws.Cells.Item(118, 7).Value = 0
$ws.Cells.Item(118, 8).Value = 10

# Row 119: San Marino
$ws.Cells.Item(119, 1).Value = "San Marino"
$ws.Cells.Item(119, 2).Value = 628
$ws.Cells.Item(119, 3).Value = 0
$ws.Cells.Item(119, 4).Value = 126
$ws.Cells.Item(119, 5).Value = 461
$ws.Cells.Item(119, 6).Value = 3
$ws.Cells.Item(119, 7).Value = 0
$ws.Cells.Item(119, 8).Value = 41
